# Completed logic of product flow report
# - Update the reporting period text
# - Insert the company name + address rows beneath the period
# - Keep the existing blank separator row (now pushed down) merged
# - Existing data rows shift down automatically as a result of the insert

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the period text. Row 2 (A2:G2) is merged in the source workbook;
#    briefly unmerge so the new text can be written, then restore the merge
#    and formatting (bordered cells, left-aligned text) across A2:G2.
$ws.Range("A2:G2").UnMerge()
$ws.Range("A2").Value = "Период: 2023-10-07 - 2023-10-29"
$ws.Range("A2:G2").Merge()
$ws.Range("A2:G2").Borders.LineStyle = 1
$ws.Range("A2:G2").HorizontalAlignment = -4131

# 2. Insert two new rows right after the period row; this pushes the old
#    blank separator row + header + data rows down by two.
$ws.Rows("3:4").Insert()

# 3. Row 3: company name, merged across A:G, same look as the period row
#    (bordered cells, left-aligned text).
$ws.Range("A3:G3").Merge()
$ws.Range("A3").Value = "ОАО Гомельский Мясокомбинат"
$ws.Range("A3:G3").Borders.LineStyle = 1
$ws.Range("A3:G3").HorizontalAlignment = -4131

# 4. Row 4: address, same formatting as row 3.
$ws.Range("A4:G4").Merge()
$ws.Range("A4").Value = "Адрес: Гомель, ул. Ильича, 2"
$ws.Range("A4:G4").Borders.LineStyle = 1
$ws.Range("A4:G4").HorizontalAlignment = -4131

# 5. Row 5 is the blank separator row that used to be row 3 before the
#    insert; it keeps its bordered-but-unaligned look and becomes merged.
$ws.Range("A5:G5").Merge()
$ws.Range("A5:G5").Borders.LineStyle = 1

Write-Output "Report header rebuilt: period/company/address rows in place."
